$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A16").Value = "Chi-Yun Chen"
$ws.Range("B16").Value = "Department of Space Science & Engineering, National Central University"
$ws.Range("C16").Value = "Taiwan astronomical Observation collaboration Platform (TOP)"

$ws.Range("B10").Select()
